$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are written with a leading apostrophe to force
# Excel to store them as text, matching the source data (avoids float
# rounding / auto-number conversion of values like "37.158.49" or "5.50").

$ws.Range("D2").Value = "'37.158.49"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "'2.061.41"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'248.47"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'0.666"
$ws.Range("E6").Value = "  +0.89%  "
$ws.Range("D7").Value = "'58.27"
$ws.Range("E7").Value = "  +4.18%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.383"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").Value = "'0.0786"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "'15.85"
$ws.Range("E12").Value = "  +1.45%  "
$ws.Range("D13").Value = "'0.911"
$ws.Range("E13").Value = "  +16.00%  "
$ws.Range("D14").Value = "'2.364.95"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'5.85"
$ws.Range("E15").Value = "  +4.83%  "
$ws.Range("D16").Value = "'2.053.25"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'18.55"
$ws.Range("E17").Value = "  +13.88%  "
$ws.Range("D18").Value = "'37.178.84"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'75.19"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("D20").Value = "'0.0₃0908"
$ws.Range("E20").Value = "  +2.04%  "
$ws.Range("D21").Value = "'5.50"
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("D22").Value = "'238.60"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  +5.42%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.66"
$ws.Range("E25").Value = "  +7.15%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "'170.98"
$ws.Range("E27").Value = "  +2.20%  "
$ws.Range("D28").Value = "'20.18"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("D29").Value = "'5.47"
$ws.Range("E29").Value = "  +18.16%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").Value = "'1.15"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").Value = "'4.88"
$ws.Range("E32").Value = "  +11.35%  "
$ws.Range("D33").Value = "'0.0626"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.0881"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'2.32"
$ws.Range("E35").Value = "  +5.62%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "'1.83"
$ws.Range("E37").Value = "  +4.54%  "
$ws.Range("D38").Value = "'1.34"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").Value = "'5.20"
$ws.Range("E39").Value = "  +6.89%  "
$ws.Range("D40").Value = "'3.12"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").Value = "'0.101"
$ws.Range("E41").Value = "  -4.14%  "
$ws.Range("D42").Value = "'0.0225"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("D43").Value = "'1.17"
$ws.Range("E43").Value = "  +5.58%  "
$ws.Range("D44").Value = "'100.13"
$ws.Range("E44").Value = "  +5.59%  "
$ws.Range("D45").Value = "'17.44"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("D46").Value = "'2.41"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'1.310.87"
$ws.Range("E47").Value = "  +2.84%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "'3.85"
$ws.Range("E48").Value = "  +15.18%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.88"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'6.95"
$ws.Range("E50").Value = "  +4.71%  "
$ws.Range("D51").Value = "'2.255.44"
$ws.Range("E51").Value = "  +1.32%  "

Write-Host "Updated cryptos list"
